# Auto-generated Excel COM-interop script applying the "New crime data collected" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Cells that change type (text "N/A" <-> number) need their format copied
#     from a same-style neighbor cell before the new value is written, so the
#     saved style index lines up with the target workbook. ---
$ws.Range("G14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = 1

$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("F14").Value = 1

$ws.Range("G14").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 1

$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("D28").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "0"

$ws.Range("D28").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = "0"

$ws.Range("H28").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"

$ws.Range("K22").Copy()
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("L22").Value = 0

$ws.Range("D26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1

$ws.Range("D15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2

$ws.Range("D15").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").Value = 2

$ws.Range("D15").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = 2

$ws.Range("D15").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$ws.Range("D15").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = 1

$ws.Range("D15").Copy()
$ws.Range("I29").PasteSpecial(-4122)
$ws.Range("I29").Value = 1

$excel.CutCopyMode = 0

# --- Plain value updates (style unchanged) ---
$ws.Range("H14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("N14").Value = -75
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = -25
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 200
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = -11.111111111111
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 17
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = -26.086956521739
$ws.Range("L16").Value = 70
$ws.Range("M16").Value = -43.333333333333
$ws.Range("N16").Value = -86.614173228346
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = -40
$ws.Range("F17").Value = 44
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = -4.347826086956
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = 5.555555555555
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 123.529411764706
$ws.Range("N17").Value = 123.529411764706
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 125
$ws.Range("I18").Value = 11
$ws.Range("J18").Value = 5
$ws.Range("K18").Value = 120
$ws.Range("L18").Value = -31.25
$ws.Range("M18").Value = -57.692307692307
$ws.Range("N18").Value = -94.179894179894
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = -63.636363636363
$ws.Range("G19").Value = 168
$ws.Range("H19").Value = -61.309523809523
$ws.Range("I19").Value = 44
$ws.Range("J19").Value = 119
$ws.Range("K19").Value = -63.025210084033
$ws.Range("L19").Value = 62.962962962963
$ws.Range("M19").Value = 22.222222222222
$ws.Range("N19").Value = -47.619047619047
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 5.263157894736
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 18
$ws.Range("K20").Value = 5.555555555555
$ws.Range("L20").Value = 216.666666666667
$ws.Range("M20").Value = 137.5
$ws.Range("N20").Value = -86.619718309859
$ws.Range("D21").Value = 62
$ws.Range("E21").Value = -30.645161290322
$ws.Range("F21").Value = 171
$ws.Range("G21").Value = 274
$ws.Range("H21").Value = -37.591240875912
$ws.Range("I21").Value = 133
$ws.Range("J21").Value = 204
$ws.Range("K21").Value = -34.803921568627
$ws.Range("L21").Value = 62.195121951219
$ws.Range("M21").Value = 11.764705882352
$ws.Range("N21").Value = -76.418439716312
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = -33.333333333333
$ws.Range("M22").Value = -50
$ws.Range("C24").Value = 52
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = 13.043478260869
$ws.Range("F24").Value = 243
$ws.Range("G24").Value = 148
$ws.Range("H24").Value = 64.189189189189
$ws.Range("I24").Value = 189
$ws.Range("J24").Value = 105
$ws.Range("K24").Value = 80
$ws.Range("L24").Value = 103.225806451613
$ws.Range("M24").Value = 119.767441860465
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 83.333333333333
$ws.Range("F25").Value = 81
$ws.Range("G25").Value = 50
$ws.Range("H25").Value = 62
$ws.Range("I25").Value = 60
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = 100
$ws.Range("L25").Value = 53.846153846153
$ws.Range("M25").Value = 76.470588235294
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 3
$ws.Range("L26").Value = -40
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 37.5
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 80
$ws.Range("N28").Value = -60
$ws.Range("N29").Value = -80
